$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.710.44"
$ws.Range("E2").Value = "  -7.29%  "

$ws.Range("D3").Value = "2.186.23"
$ws.Range("E3").Value = "  -7.36%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.44"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("E6").Value = "  -7.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.42"
$ws.Range("E7").Value = "  -7.85%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.37"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").Value = "  -7.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.43"
$ws.Range("E12").Value = "  -5.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -5.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.55"
$ws.Range("E14").Value = "  -10.02%  "

$ws.Range("D15").Value = "2.511.67"
$ws.Range("E15").Value = "  -7.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.57"
$ws.Range("E16").Value = "  -10.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.827"
$ws.Range("E17").Value = "  -9.90%  "

$ws.Range("D18").Value = "2.186.87"
$ws.Range("E18").Value = "  -7.55%  "

$ws.Range("D19").Value = "40.684.41"

$ws.Range("E20").Value = "  -9.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.41"
$ws.Range("E21").Value = "  -7.14%  "

$ws.Range("E22").Value = "  -8.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.01"
$ws.Range("E23").Value = "  -9.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +6.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("E26").Value = "  -5.08%  "

$ws.Range("E27").Value = "  -4.08%  "

$ws.Range("E28").Value = "  -5.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.64"
$ws.Range("E29").Value = "  -8.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.46"
$ws.Range("E30").Value = "  -4.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.13"
$ws.Range("E31").Value = "  -9.89%  "

$ws.Range("E32").Value = "  -10.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  -8.36%  "

$ws.Range("E34").Value = "  -7.33%  "

$ws.Range("E35").Value = "  -5.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.55"
$ws.Range("E36").Value = "  -10.59%  "

$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.33"
$ws.Range("E38").Value = "  +14.73%  "

$ws.Range("E39").Value = "  -6.70%  "

$ws.Range("E40").Value = "  -4.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("E41").Value = "  -12.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "62.24"
$ws.Range("E42").Value = "  -4.44%  "

$ws.Range("E43").Value = "  -11.58%  "

$ws.Range("E44").Value = "  -5.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.189"
$ws.Range("E45").Value = "  -6.94%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0978"
$ws.Range("E47").Value = "  -8.21%  "

$ws.Range("E48").Value = "  +1.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.23"
$ws.Range("E49").Value = "  +6.04%  "

$ws.Range("E50").Value = "  -6.84%  "

$ws.Range("E51").Value = "  -6.44%  "

